$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new header columns (AC, AD, AE) reusing the same header
# formatting (bold, centered, bordered) as the existing header cells by
# copying an existing header cell's formatting onto the new range first.
$ws.Range("AB1").Copy($ws.Range("AC1:AE1"))

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row.
$firstRow = 2
$lastRow = 45

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $ws.Range("AC$row").Value = 64
    $ws.Range("AD$row").Value = 98
    $ws.Range("AE$row").Value = 0
}
